$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "76.332.82"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.034.38"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +4.26%  "

$ws.Range("E4").Value = "  +0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "197.77"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "616.82"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.43%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.88%  "

$ws.Range("E9").Value = "  +5.16%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.033.39"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.35%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.432"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("E12").Value = "  -0.71%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.23"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.93%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.594.49"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.39%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "28.65"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.59%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "75.880.06"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000191"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.09%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.034.04"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +4.31%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.50"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.14%  "

$ws.Range("E20").Value = "  +2.93%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "378.10"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.91%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.32"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.11%  "

$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("E24").Value = "  +3.92%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "72.56"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("E26").Value = "  -0.08%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "4.31"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.66"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.0000107"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("E30").Value = "  +0.05%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.22"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.06%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.38"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "490.31"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.38%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.91"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.52%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "20.48"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.97%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "162.95"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("E38").Value = "  +2.12%  "

$ws.Range("E39").Value = "  +5.21%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.379"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.14%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "190.43"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.63%  "

$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("E43").Value = "  +0.00%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.794"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +21.52%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "5.04"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "41.67"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.84%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.36%  "

$ws.Range("E48").Value = "  -1.18%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.38"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.44%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.596"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.41%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.86"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.83%  "
